# "create db at runtime"
# Update the status of issue #8 (row 9) from "Re-open" to "Fixed",
# and leave the active selection on that cell (F9).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F9").Value = "Fixed"

$ws.Activate()
$ws.Range("F9").Select()
